$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.465.80"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "3.642.41"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "196.85"
$ws.Range("E5").Value = "  +6.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "581.46"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("D7").Value = "3.637.03"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.681"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  +6.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.47"
$ws.Range("E12").Value = "  +4.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000293"
$ws.Range("E13").Value = "  +14.94%  "
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "4.228.00"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "3.646.40"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.61"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "68.434.26"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.68"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.56"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.99"
$ws.Range("E23").Value = "  +25.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.26"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.19"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.97"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.65"
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("E28").Value = "  +6.86%  "
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("E30").Value = "  +18.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.20"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "701.81"
$ws.Range("E33").Value = "  +17.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.27"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.86"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.85"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("E38").Value = "  +12.84%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "0.0₃0789"
$ws.Range("E40").Value = "  +6.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  +18.59%  "
$ws.Range("E42").Value = "  +4.10%  "
$ws.Range("D43").Value = "3.228.51"
$ws.Range("E43").Value = "  +19.44%  "
$ws.Range("E44").Value = "  +12.85%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.00"
$ws.Range("E46").Value = "  +32.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0423"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.93"
$ws.Range("E48").Value = "  +7.86%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.133"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.46"
$ws.Range("E51").Value = "  +3.98%  "
